$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TruthFulQA Self-refine")
$tbl = $ws.ListObjects.Item(1)

# --- 1. Remove the two rows whose data is dropped entirely in the new table ---
# Old row 3  = "History_3"   (duplicate of History_3_old, dropped)
# Old row 10 = "New_4_old"   (dropped)
# Delete row 10 first (higher index) so row 3's index doesn't shift before we delete it.
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(3).Delete()

# After the two deletions the 10 remaining data rows (rows 2-11) are, in order:
#  2 History_2       3 History_3_old   4 History_4   5 Hisory
#  6 New_2           7 New_3           8 New_4       9 New_5
# 10 New             11 Original
# We now overwrite every cell (A:E) of rows 1-11 with the final target content,
# in the new row order: header, 0,1,2,3,4,5,H1,H2,H3,H4

# --- 2. Header row ---
$ws.Cells.Item(1,1).Value = "version"
$ws.Cells.Item(1,2).Value = "bleu_acc"
$ws.Cells.Item(1,3).Value = "rouge1_acc"
$ws.Cells.Item(1,4).Value = "bleurt_acc"
$ws.Cells.Item(1,5).Value = "Average"

# --- 3. Data rows: version, bleu_acc, rouge1_acc, bleurt_acc ---
$data = @(
  @("0", 0.43574051407588699, 0.438188494492044,   0.58506731946144397),
  @("1", 0.38800489596083199, 0.428396572827417,   0.48714810281517701),
  @("2", 0.446756425948592,   0.46878824969400201, 0.53121175030599699),
  @("3", 0.50795593635250902, 0.53488372093023195, 0.61689106487148104),
  @("4", 0.494492044063647,   0.52509179926560501, 0.581395348837209),
  @("5", 0.438188494492044,   0.45042839657282702, 0.52509179926560501),
  @("H1",0.26805385556915501, 0.42105263157894701, 0.52753977968176202),
  @("H2",0.69889840881272902, 0.75764993880048903, 0.77600979192166397),
  @("H3",0.70379436964504205, 0.77233782129742901, 0.78212974296205595),
  @("H4",0.61933904528763695, 0.70991432068543403, 0.73929008567931398)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $row = $data[$i]
  $ws.Cells.Item($r,1).Value = $row[0]
  $ws.Cells.Item($r,2).Value = $row[1]
  $ws.Cells.Item($r,3).Value = $row[2]
  $ws.Cells.Item($r,4).Value = $row[3]
  $ws.Cells.Item($r,5).Formula = "=AVERAGE(Table4[[#This Row],[bleu_acc]:[bleurt_acc]])"
}

# --- 4. Formatting: version column centered, numeric columns to 0.0000 ---
$tbl.ListColumns.Item(1).DataBodyRange.HorizontalAlignment = -4108
$ws.Range("B2:D11").NumberFormat = "0.0000"
$ws.Range("E2:E11").NumberFormat = "0.0000"

# --- 5. Sheet view: zoom + selection ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$ws.Range("A8:E8").Select()

# --- 6. Column widths ---
$ws.Columns.Item(1).ColumnWidth = 6.166666666666667
$ws.Columns.Item(2).ColumnWidth = 8.833333333333334
$ws.Columns.Item(3).ColumnWidth = 10.333333333333334
$ws.Columns.Item(4).ColumnWidth = 9.666666666666666

$wb.Save()
